# Updated symbol list on Fri Dec 30 19:56:18 UTC 2022 with GitHub Actions
#
# Re-applies the latest scrape of coinranking.com data onto the "cryptos"
# worksheet: refreshed Price (column D) figures, plus the #10/#19 ranking
# swap between MandalaExchangeToken and LiechtensteinCryptoassetsExchange
# (including their Coin/Link/Volume(1h) cells), and a couple of
# "...Bestin24h" / "...Worstin24h" suffix changes in Volume(1h).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Cell,
        [string]$Value
    )
    # Prefix with an apostrophe so numeric-looking strings (prices such as
    # "244.96") are stored as text, not coerced to a Number -- the source
    # sheet keeps every value (even prices) as plain text cells.
    $range = $ws.Range($Cell)
    $range.Value = "'" + $Value
    # Writing a text-quoted value stamps the cell with a "@" (text) number
    # format / style; reset back to the workbook's default "Normal" style
    # so the cell matches its unedited neighbours (no stray style index).
    $range.Style = "Normal"
}

# Row 2 - BNB
Set-TextCell "D2" "244.96"

# Row 4 - HuobiToken
Set-TextCell "D4" "4.979"

# Row 6 - KuCoinToken
Set-TextCell "D6" "6.545"

# Row 7 - GateToken
Set-TextCell "D7" "3.005"

# Row 8 - MXToken
Set-TextCell "D8" "0.8112"

# Row 9 - FTXToken
Set-TextCell "D9" "0.8400"

# Row 10 - WazirX
Set-TextCell "D10" "0.1337"

# Row 11 - was MandalaExchangeToken, now LiechtensteinCryptoassetsExchange
Set-TextCell "B11" "LiechtensteinCryptoassetsExchange"
Set-TextCell "C11" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D11" "0.03242"
Set-TextCell "E11" "10LiechtensteinCryptoassetsExchangeLCX"

# Row 12 - BitrueCoin
Set-TextCell "D12" "0.02848"

# Row 13 - BitMartToken
Set-TextCell "D13" "0.09410"

# Row 14 - BitForexToken
Set-TextCell "D14" "0.001526"

# Row 15 - One
Set-TextCell "D15" "0.0005989"

# Row 16 - TigerCash
Set-TextCell "D16" "0.006226"

# Row 17 - LEO
Set-TextCell "D17" "3.499"

# Row 19 - BitpandaEcosystemToken
Set-TextCell "D19" "0.3198"

# Row 20 - was LiechtensteinCryptoassetsExchange, now MandalaExchangeToken
Set-TextCell "B20" "MandalaExchangeToken"
Set-TextCell "C20" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D20" "0.06966"
Set-TextCell "E20" "19MandalaExchangeTokenMDX"

# Row 22 - MCDex
Set-TextCell "D22" "3.738"

# Row 23 - CoinExToken
Set-TextCell "D23" "0.04694"

# Row 25 - BitKan
Set-TextCell "D25" "0.001242"

# Row 26 - HotbitToken
Set-TextCell "D26" "0.004527"

# Row 27 - NitroEx
Set-TextCell "D27" "0.00009698"
Set-TextCell "E27" "26NitroExNTX"

# Row 28 - UpBots
Set-TextCell "D28" "0.0001940"

# Row 40 - IDEX
Set-TextCell "D40" "0.03648"

# Row 41 - BKEXToken
Set-TextCell "D41" "0.1352"

# Row 42 - KickToken
Set-TextCell "D42" "0.006229"
Set-TextCell "E42" "41KickTokenKICKBestin24h"

# Row 44 - LocalTraders
Set-TextCell "D44" "0.008074"

# Row 45 - CoinLion
Set-TextCell "D45" "0.00005275"

# Row 49 - CryptobidCoin
Set-TextCell "D49" "0.00002100"

# Row 50 - SpecialPowerGold
Set-TextCell "D50" "0.0002000"
